# Auto-generated edit script: updates market-price columns (H-N) across all 8 sheets
# to match the refreshed data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3851128.8
$ws.Range("I43").Value = 5131505
$ws.Range("K43").Value = 5131505
$ws.Range("M43").Value = -5131436
$ws.Range("H70").Value = 5444
$ws.Range("J70").Value = 5650
$ws.Range("L70").Value = 16950
$ws.Range("N70").Value = -17490
$ws.Range("H73").Value = 5444
$ws.Range("J73").Value = 5650
$ws.Range("L73").Value = 16950
$ws.Range("N73").Value = -18822
$ws.Range("H111").Value = 1502
$ws.Range("I111").Value = 1379.4445
$ws.Range("J111").Value = 1659.5714
$ws.Range("K111").Value = 4138.333500000001
$ws.Range("L111").Value = 4978.7142
$ws.Range("M111").Value = -1071.333500000001
$ws.Range("N111").Value = -11112.7142

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3596.9734
$ws.Range("I32").Value = 1830.1562
$ws.Range("K32").Value = 1830.1562
$ws.Range("M32").Value = -1543.1562
$ws.Range("H45").Value = 3087.5
$ws.Range("I45").Value = 2796.5715
$ws.Range("K45").Value = 2796.5715
$ws.Range("M45").Value = -2419.5715
$ws.Range("H61").Value = 36202
$ws.Range("I61").Value = 44127.75
$ws.Range("K61").Value = 44127.75
$ws.Range("M61").Value = -43915.75
$ws.Range("H102").Value = 286452.44
$ws.Range("I102").Value = 473147
$ws.Range("K102").Value = 473147
$ws.Range("M102").Value = -471525
$ws.Range("H106").Value = 83000
$ws.Range("J106").Value = 83000
$ws.Range("L106").Value = 83000
$ws.Range("N106").Value = -85524
$ws.Range("H122").Value = 2203.6316
$ws.Range("I122").Value = 2261.8125
$ws.Range("K122").Value = 6785.4375
$ws.Range("M122").Value = -4335.4375
$ws.Range("H132").Value = 73418.44500000001
$ws.Range("I132").Value = 142891.75
$ws.Range("J132").Value = 17839.8
$ws.Range("K132").Value = 428675.25
$ws.Range("L132").Value = 53519.39999999999
$ws.Range("M132").Value = -426145.25
$ws.Range("N132").Value = -58579.39999999999
$ws.Range("H136").Value = 36202
$ws.Range("I136").Value = 44127.75
$ws.Range("K136").Value = 132383.25
$ws.Range("M136").Value = -129833.25

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2564.4
$ws.Range("I20").Value = 3069
$ws.Range("J20").Value = 1176.75
$ws.Range("K20").Value = 3069
$ws.Range("L20").Value = 1176.75
$ws.Range("M20").Value = -2822
$ws.Range("N20").Value = -1670.75
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").Value = ""
$ws.Range("H94").Value = 1523054.4
$ws.Range("I94").Value = 2283869
$ws.Range("J94").Value = 1425
$ws.Range("K94").Value = 2283869
$ws.Range("L94").Value = 1425
$ws.Range("M94").Value = -2283418
$ws.Range("N94").Value = -2327
$ws.Range("H107").Value = 1350.6875
$ws.Range("I107").Value = 1250.76
$ws.Range("J107").Value = 1707.5714
$ws.Range("K107").Value = 1250.76
$ws.Range("L107").Value = 1707.5714
$ws.Range("M107").Value = 669.24
$ws.Range("N107").Value = -5547.5714
$ws.Range("H134").Value = 4112
$ws.Range("I134").Value = 2116.8572
$ws.Range("K134").Value = 6350.571599999999
$ws.Range("M134").Value = -3815.571599999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2648.25
$ws.Range("I134").Value = 2776.3555
$ws.Range("J134").Value = 1824.7142
$ws.Range("K134").Value = 8329.066500000001
$ws.Range("L134").Value = 5474.142599999999
$ws.Range("M134").Value = -5794.066500000001
$ws.Range("N134").Value = -10544.1426

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 844.6667
$ws.Range("I5").Value = 742
$ws.Range("J5").Value = 1050
$ws.Range("K5").Value = 2226
$ws.Range("L5").Value = 3150
$ws.Range("M5").Value = -2114
$ws.Range("N5").Value = -3374
$ws.Range("H70").Value = 6360.8
$ws.Range("I70").Value = 4153
$ws.Range("K70").Value = 12459
$ws.Range("M70").Value = -12144
$ws.Range("H73").Value = 6360.8
$ws.Range("I73").Value = 4153
$ws.Range("K73").Value = 12459
$ws.Range("M73").Value = -11367
$ws.Range("H113").Value = 1010
$ws.Range("J113").Value = 1222
$ws.Range("L113").Value = 3666
$ws.Range("N113").Value = -8006
$ws.Range("H119").Value = 17399.25
$ws.Range("I119").Value = 14800
$ws.Range("J119").Value = 19998.5
$ws.Range("K119").Value = 44400
$ws.Range("L119").Value = 59995.5
$ws.Range("M119").Value = -39562
$ws.Range("N119").Value = -69671.5
$ws.Range("H121").Value = 751107.75
$ws.Range("I121").Value = 1599
$ws.Range("J121").Value = 1500616.5
$ws.Range("K121").Value = 4797
$ws.Range("L121").Value = 4501849.5
$ws.Range("M121").Value = -3487
$ws.Range("N121").Value = -4504469.5
$ws.Range("H135").Value = 844.6667
$ws.Range("I135").Value = 742
$ws.Range("J135").Value = 1050
$ws.Range("K135").Value = 6678
$ws.Range("L135").Value = 9450
$ws.Range("M135").Value = -4143
$ws.Range("N135").Value = -14520

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = ""
$ws.Range("H63").Value = 99999
$ws.Range("J63").Value = 99999
$ws.Range("L63").Value = 99999
$ws.Range("N63").Value = -101371
$ws.Range("H66").Value = 99999
$ws.Range("J66").Value = 99999
$ws.Range("L66").Value = 299997
$ws.Range("N66").Value = -306861
$ws.Range("H69").Value = 22000
$ws.Range("I69").Value = 22000
$ws.Range("K69").Value = 22000
$ws.Range("M69").Value = -21251
$ws.Range("H72").Value = 22000
$ws.Range("I72").Value = 22000
$ws.Range("K72").Value = 66000
$ws.Range("M72").Value = -62256
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""
$ws.Range("H93").Value = 44147.855
$ws.Range("J93").Value = 44147.855
$ws.Range("L93").Value = 44147.855
$ws.Range("N93").Value = -47891.855
$ws.Range("H94").Value = 40112
$ws.Range("J94").Value = 31816
$ws.Range("L94").Value = 31816
$ws.Range("N94").Value = -33168
$ws.Range("H132").Value = 10016.077
$ws.Range("I132").Value = 10301.667
$ws.Range("J132").Value = 9771.286
$ws.Range("K132").Value = 30905.001
$ws.Range("L132").Value = 29313.858
$ws.Range("M132").Value = -28375.001
$ws.Range("N132").Value = -34373.858

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8001880
$ws.Range("I16").Value = 9092462
$ws.Range("J16").Value = 4276.3335
$ws.Range("K16").Value = 9092462
$ws.Range("L16").Value = 4276.3335
$ws.Range("M16").Value = -9092292
$ws.Range("N16").Value = -4616.3335
$ws.Range("H122").Value = 49692520
$ws.Range("I122").Value = 52634308
$ws.Range("K122").Value = 157902924
$ws.Range("M122").Value = -157900474
$ws.Range("H136").Value = 2416.3215
$ws.Range("I136").Value = 1871.4231
$ws.Range("K136").Value = 5614.2693
$ws.Range("M136").Value = -3064.2693

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = ""
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = ""
$ws.Range("H81").Value = 2983775.8
$ws.Range("I81").Value = 1901848.1
$ws.Range("J81").Value = 6950844
$ws.Range("K81").Value = 3803696.2
$ws.Range("L81").Value = 13901688
$ws.Range("M81").Value = -3802635.2
$ws.Range("N81").Value = -13903810
$ws.Range("H84").Value = 2983775.8
$ws.Range("I84").Value = 1901848.1
$ws.Range("J84").Value = 6950844
$ws.Range("K84").Value = 19018481
$ws.Range("L84").Value = 69508440
$ws.Range("M84").Value = -19013177
$ws.Range("N84").Value = -69519048

